$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("BasicSearch")
$ws2 = $wb.Worksheets.Item("AdvanceSearch")

# --- AdvanceSearch (sheet2): fix the ALL (Acute Lymphoblastic Leukemia) row's
# CancerStageId / CancerStage first ---
$ws2.Range("E4").Value = "C142810"
$ws2.Range("F4").Value = "Recurrent Acute Lymphoblastic Leukemia"

# --- BasicSearch (sheet1): Age / ZipCode become text (quote-prefixed) values ---
$ws1.Range("B2").Value = "'65"
$ws1.Range("B2").NumberFormat = "@"
$ws1.Range("C2").Value = "'20105"
$ws1.Range("C2").NumberFormat = "@"
$ws1.Range("B3").Value = "'50"
$ws1.Range("B3").NumberFormat = "@"
$ws1.Range("C3").Value = "'20105"
$ws1.Range("C3").NumberFormat = "@"

# --- AdvanceSearch (sheet2): CancerSubTypeId / CancerSubType for the ALL row ---
$ws2.Range("C4").Value = "C122614"
$ws2.Range("D4").Value = "Infant Acute Lymphoblastic Leukemia"

# --- BasicSearch (sheet1): add a new "Lung Cancer" row ---
$ws1.Range("A4").Value = "Lung Cancer"
$ws1.Range("B4").Value = "'120"
$ws1.Range("C4").Value = "'20105"

# --- AdvanceSearch (sheet2): new Country / State / City / StateCode columns ---
$ws2.Range("J1").Value = "Country"
$ws2.Range("K1").Value = "State"
$ws2.Range("M1").Value = "City"
$ws2.Range("J2:J4").Value = "United States"
$ws2.Range("K2:K4").Value = "California"
$ws2.Range("M2:M4").Value = "Los Angeles"
$ws2.Range("L1").Value = "StateCode"
$ws2.Range("L2:L4").Value = "CA"

# --- AdvanceSearch (sheet2): new Hospital column ---
$ws2.Range("N1").Value = "Hospital"
$ws2.Range("N2:N4").Value = "John F Kennedy Medical Center"

# --- AdvanceSearch (sheet2): new Drug / DrugId columns ---
$ws2.Range("O1").Value = "Drug"
$ws2.Range("O2:O4").Value = "Antineoplastic Agent"
$ws2.Range("P1").Value = "DrugId"
$ws2.Range("P2:P4").Value = "C274"

# --- AdvanceSearch (sheet2): new Treatment / TreatmentId columns (per-row values) ---
$ws2.Range("Q2").Value = "Radiation Therapy"
$ws2.Range("R2").Value = "C15313"
$ws2.Range("Q1").Value = "Treatment"
$ws2.Range("R1").Value = "TreatmentId"
$ws2.Range("Q3").Value = "Physical Therapy"
$ws2.Range("R3").Value = "C15302"
$ws2.Range("R4").Value = "C94626"
$ws2.Range("Q4").Value = "Chemoradiotherapy"

# --- AdvanceSearch (sheet2): new TrialPhase column ---
$ws2.Range("S1").Value = "TrialPhase"
$ws2.Range("S2").Value = "Phase I"
$ws2.Range("S3").Value = "Phase II"
$ws2.Range("S4").Value = "Phase III"

# NOTE: H1/I1, G2/H2/I2, G3/H3/I3, A4/B4/G4/H4/I4 keep their original text —
# only their shared-string index shifts once unused strings are pruned on
# save, so no explicit re-write is required for those cells.

# --- View/selection bookkeeping to match the saved workbook state ---
$ws1.Range("C11").Select()
$ws2.Activate()
$ws2.Range("N10").Select()
